$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: stamp the date-column (A) style onto the new rows (50-68) by copying
# the existing styled cell A2 (bold/centered/bordered) down to those rows before
# we overwrite their values, so the new cells inherit the same cell style (s="1").
for ($r = 50; $r -le 68; $r++) {
    $ws.Range("A2").Copy($ws.Range("A" + $r))
    $ws.Range("B2").Copy($ws.Range("B" + $r))
}

# Step 2: write the final data set (rows 2-68) in the order required by the target
# layout -- for each year, the Oct/Nov/Dec rows now lead (inserted ahead of Jan..Sep),
# and 2023 Jan-Jul are appended as brand-new rows at the end.
$data = @(
    @("2018-10", 113.6, 116.9, 126.5),
    @("2018-11", 123.7, 102.6, 116.4),
    @("2018-12", 113.4, 98.7, 104.2),
    @("2018-01", 118.0131, 97.9347, 109.1769),
    @("2018-02", 114.4, 103.1, 109),
    @("2018-03", 114.6, 106.4, 107.1),
    @("2018-04", 107.8, 103.4, 110.6),
    @("2018-05", 108.6, 107.9, 115.9),
    @("2018-06", 117.7, 107.4, 120.2),
    @("2018-07", 118.5, 115.1, 125.9),
    @("2018-08", 116.2, 117.7, 124.2),
    @("2018-09", 114.4, 115.5, 126.5),
    @("2019-10", 88.7, 95.2, 87.9),
    @("2019-11", 84.3, 101.7, 91.3),
    @("2019-12", 87, 106.3, 100.5),
    @("2019-01", 105.9, 106.1, 96.9),
    @("2019-02", 107.3, 100.9, 97.6),
    @("2019-03", 107.5, 102.8, 102.7),
    @("2019-04", 108.1, 105.2, 103.5),
    @("2019-05", 109.1, 101, 101.2),
    @("2019-06", 105.6, 95.7, 96.7),
    @("2019-07", 100.3, 98.1, 93.9),
    @("2019-08", 98.6, 93.6, 93.3),
    @("2019-09", 90.9, 97.2, 90.3),
    @("2020-10", 100, 82.8, 78),
    @("2020-11", 107.7, 82.2, 79.2),
    @("2020-12", 112.5, 82, 82.9),
    @("2020-01", 91.1, 105.5, 107.6),
    @("2020-02", 91.9, 113.7, 100.7),
    @("2020-03", 90.3, 99.1, 89.2),
    @("2020-04", 88.8, 92.6, 78.6),
    @("2020-05", 86.1, 87.4, 73.6),
    @("2020-06", 87.5, 91.1, 76.6),
    @("2020-07", 92.6, 86.1, 81.1),
    @("2020-08", 91, 88.9, 81.8),
    @("2020-09", 95.4, 87.3, 80.7),
    @("2021-10", 195.3, 113, 144),
    @("2021-11", 169.1, 115.9, 149.7),
    @("2021-12", 137, 109.8, 136.4),
    @("2021-01", 123.6, 86.5, 85.6),
    @("2021-02", 131, 91.1, 94),
    @("2021-03", 128.7, 101.5, 111.2),
    @("2021-04", 124, 110.7, 123.8),
    @("2021-05", 144.3, 118.2, 132.4),
    @("2021-06", 145, 116.8, 134.4),
    @("2021-07", 143.3, 112.7, 133.9),
    @("2021-08", 154.4, 110.6, 131.6),
    @("2021-09", 179.6, 111.4, 132.6),
    @("2022-10", 74.7, 106.4, 117.9),
    @("2022-11", 78.2, 101.9, 113.6),
    @("2022-12", 95.1, 101.5, 113),
    @("2022-01", 128.8, 111.4, 130.6),
    @("2022-02", 118.4, 110, 132.9),
    @("2022-03", 133.6, 111.6, 132.5),
    @("2022-04", 159.3, 110.7, 134.9),
    @("2022-05", 139.7, 110, 132.9),
    @("2022-06", 125.1, 111.1, 136.5),
    @("2022-07", 113.2, 110.7, 131.4),
    @("2022-08", 95.4, 107.9, 126.4),
    @("2022-09", 80.6, 107.6, 126.5),
    @("2023-01", 94.7, 98.5, 108.3),
    @("2023-02", 94.5, 99, 102.9),
    @("2023-03", 87.5, 96.5, 94.9),
    @("2023-04", 75.4, 94, 91.1),
    @("2023-05", 64.7, 93.3, 89),
    @("2023-06", 63.3, 90.9, 83),
    @("2023-07", 68.6, 91.8, 84.1)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
    $row++
}

Write-Output ("Dimension rows: " + $ws.UsedRange.Rows.Count)
